# Updated cryptos list values (prices and 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.904.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.888.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7330'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3105'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06900'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7656'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.899.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.232'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.920.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.746'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007759'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.154.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.914'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.303'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1270'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.016'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.21%  '
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.082'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05088'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.276'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7365'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.721'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01924'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.774'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.326'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4449'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.930'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8365'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.600'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.771'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.034.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '943.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.71%  '
